# Fix typos/accents in the "Journal" log (Description column) and the
# "Descritpion" -> "Description" header, matching the "updated for last
# version" commit. Edits are applied in the same order they were made in
# the original commit so the shared-string table comes out in the same
# order Excel would produce.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Mise en place d'Eedomus"
$ws.Range("A1").Value = "Description"
$ws.Range("A2").Value = "Installation du Raspberry"
$ws.Range("A5").Value = "Réinstallation Raspberry pi et domoticz"
$ws.Range("A6").Value = "Théorie câblage plus début de la maquette"
$ws.Range("A7").Value = "Maquette presque terminée, circuit 12v restant"
$ws.Range("A8").Value = "Fin de la maquette sans qubino, début de l'ajout du qubino"
$ws.Range("A10").Value = "Ajout des appareils enocean et début des scénarios"
$ws.Range("A15").Value = "Finalisation de la documentation"

$ws.Columns.Item(1).ColumnWidth = 58.5

$ws.Range("F18").Select()
